$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("G8").Value = 3
$ws.Range("I8").Value = 2.3
$ws.Range("J8").Value = 3.75
$ws.Range("L8").Value = 3.1
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 9
$ws.Range("Q8").Value = 2.15
$ws.Range("R8").Value = 1.67
$ws.Range("S8").Value = 1.44
$ws.Range("T8").Value = 2.63
$ws.Range("W8").Value = 9
$ws.Range("X8").Value = 15
$ws.Range("Y8").Value = 12
$ws.Range("Z8").Value = 34
$ws.Range("AB8").Value = 41
$ws.Range("AI8").Value = 9.5
$ws.Range("AJ8").Value = 21
$ws.Range("AN8").Value = 5
$ws.Range("AT8").Value = 2.63
$ws.Range("AU8").Value = 8
$ws.Range("AW8").Value = 4.33
$ws.Range("AX8").Value = 13
$ws.Range("AY8").Value = 23
$ws.Range("AZ8").Value = 41

# Row 11
$ws.Range("H11").Value = 3.1
$ws.Range("I11").Value = 4
$ws.Range("K11").Value = 2.05
$ws.Range("M11").Value = 1.08
$ws.Range("N11").Value = 8
$ws.Range("O11").Value = 1.36
$ws.Range("P11").Value = 3
$ws.Range("Q11").Value = 2.2
$ws.Range("R11").Value = 1.65
$ws.Range("U11").Value = 1.91
$ws.Range("V11").Value = 1.8
$ws.Range("W11").Value = 6.5
$ws.Range("X11").Value = 9
$ws.Range("AA11").Value = 19
$ws.Range("AB11").Value = 34
$ws.Range("AC11").Value = 8
$ws.Range("AG11").Value = 10
$ws.Range("AI11").Value = 15
$ws.Range("AM11").Value = 351
$ws.Range("AR11").Value = 67
$ws.Range("AS11").Value = 201
$ws.Range("AU11").Value = 8.5
$ws.Range("AV11").Value = 67
$ws.Range("AY11").Value = 34
$ws.Range("BB11").Value = 251

# Row 12
$ws.Range("G12").Value = 1.4
$ws.Range("H12").Value = 4.2
$ws.Range("I12").Value = 9
$ws.Range("J12").Value = 1.95
$ws.Range("Q12").Value = 2.1
$ws.Range("R12").Value = 1.7
$ws.Range("AA12").Value = 15
$ws.Range("AC12").Value = 8
$ws.Range("AG12").Value = 17
$ws.Range("AI12").Value = 26
$ws.Range("AJ12").Value = 101
$ws.Range("AK12").Value = 67
$ws.Range("AO12").Value = 7
$ws.Range("BA12").Value = 251
